$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.716.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.59%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.074.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.12%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'233.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.74%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'58.62"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.43%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0783"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.39%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +3.39%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.380.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'20.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.98%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.26%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.73%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.078.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.94%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.654.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.58%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.05%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'71.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.07%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0834"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.22%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'228.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.70%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.62%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'171.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Cosmos"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.81%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Kaspa"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'0.137"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.57%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.53%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -2.97%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.41%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.20%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.52%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -5.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.38%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.16%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.00%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'99.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.93%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0972"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.74%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.96%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.45%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'16.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.49%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.441.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.67%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'4.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.58%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.45%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'3.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.72%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.265.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.16%  "
$ws.Range("E51").Style = "Normal"
